$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.193.69'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +7.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.569.05'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.98%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '416.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.15'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.650'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +4.13%  '
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.556.09'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.779'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.180'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +28.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000332'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +52.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.62'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.89'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.135.08'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.18%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.18'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.545.14'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.12'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.90%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.152.50'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.24%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.40'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '462.16'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '90.14'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.17'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.96'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.36'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.87'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.05%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.93'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.67%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.31'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.116'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.27'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.19%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.07'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.49'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.76%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0491'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0774'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +37.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.147'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +9.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '147.99'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.72'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.96'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.32'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.21'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.306'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.96'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '121.39'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +10.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.26'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.58'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +10.02%  '
